# "Add files via upload" -- refresh the "Historias de Usuarios" (HU) sheet:
#   - HU-01..HU-05/HU-06/HU-07 wording: "deseo" -> "quiero"
#   - insert a new yellow separator row (row 8)
#   - append three new AUXILIAR... wait, PACIENTE user stories (HU-08/09/10)
#   - append a trailing yellow separator row (row 12) and select it

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rows 1-7: swap "deseo" -> "quiero" in column B, keep the bold "AUXILIAR"
#    run (characters 6-13) that Excel stores as rich text.
# ---------------------------------------------------------------------------
$bText = @{
    1 = "Como AUXILIAR del hospital, quiero registrar la informacion de un nuevo paciente para que quede oficialmente registrado en el programa de hospitalizacion en casa"
    2 = "Como AUXILIAR del hospital, quiero registrar la informacion de un nuevo familiar designado para que quede oficialmente registrado en el programa de hospitalizacion en casa"
    3 = "Como AUXILIAR del hospital, quiero registrar la informacion de un nuevo medico para que quede oficialmente registrado en el programa de hospitalizacion en casa"
    4 = "Como AUXILIAR del hospital, quiero registrar la informacion de un nuevo enfermero para que quede oficialmente registrado en el programa de hospitalizacion en casa"
    5 = "Como AUXILIAR del hospital, quiero consultar la informacion de un paciente para verificar los datos registrados en el programa hospitalizacion en casa y tener sus datos de contacto en caso de necesidad"
    6 = "Como AUXILIAR del hospital, quiero asignar un enfermero a un paciente para que tenga una persona de la salud para su seguimiento"
    7 = "Como AUXILIAR del hospital, quiero asignar un medico a un paciente para que tenga un responsable de sus salud en el programa hospitalizacion en casa"
}

foreach ($r in 1..7) {
    $cell = $ws.Range("B$r")
    $txt = $bText[$r]
    $cell.Value = $txt
    # three runs: "Como " / bold role-name (8 chars, starting at char 6) / the rest
    $cell.Characters(6, 8).Font.Bold = $true
    $cell.Characters(1, 5).Font.Size = 11
    $cell.Characters(14, $txt.Length - 13).Font.Size = 11
}

# ---------------------------------------------------------------------------
# 2) Row 8 becomes a short, yellow-filled blank separator row.
# ---------------------------------------------------------------------------
$sep1 = $ws.Range("A8:B8")
$sep1.Interior.Color = 65535
$sep1.NumberFormat = "0#"
$sep1.Font.Bold = $true
$sep1.HorizontalAlignment = -4108
$sep1.VerticalAlignment = -4108
$ws.Rows(8).RowHeight = 15

# ---------------------------------------------------------------------------
# 3) Rows 9-11: new HU-08 / HU-09 / HU-10 (PACIENTE) user stories.
# ---------------------------------------------------------------------------
$aText = @{
    9  = "HU-08"
    10 = "HU-09"
    11 = "HU-10"
}
$bText2 = @{
    9  = "Como PACIENTE del programa Hospitalizacion en Casa, quiero registrar mis signos vitales para que sirvan como referencia para mi cuidado en casa"
    10 = "Como PACIENTE del programa Hospitalizacion en Casa, quiero actualizar mis datos personales para mantener al dia mis datos demograficos en caso de ser necesario."
    11 = "Como PACIENTE del programa Hospitalizacion en Casa, quiero consultar mi historia clinica para verificar los datos registrados en el programa."
}

foreach ($r in 9..11) {
    $ws.Range("A$r").Value = $aText[$r]
}
foreach ($r in 9..11) {
    $cell = $ws.Range("B$r")
    $txt = $bText2[$r]
    $cell.Value = $txt
    $cell.Characters(6, 8).Font.Bold = $true
    $cell.Characters(1, 5).Font.Size = 11
    $cell.Characters(14, $txt.Length - 13).Font.Size = 11
}

# Row 9/10 already carried the HU-formatted style from the template, but row
# 11 is brand new -- give it the same look as the rows above it.
$a11 = $ws.Range("A11")
$a11.NumberFormat = "0#"
$a11.Font.Bold = $true
$a11.HorizontalAlignment = -4108
$a11.VerticalAlignment = -4108

$b11 = $ws.Range("B11")
$b11.HorizontalAlignment = -4130
$b11.VerticalAlignment = -4130

foreach ($r in 9..11) {
    $ws.Rows($r).RowHeight = 30
}

# ---------------------------------------------------------------------------
# 4) Row 12: trailing yellow separator row, then leave it selected.
# ---------------------------------------------------------------------------
$sep2 = $ws.Range("A12:B12")
$sep2.Interior.Color = 65535
$sep2.NumberFormat = "0#"
$sep2.Font.Bold = $true
$sep2.HorizontalAlignment = -4108
$sep2.VerticalAlignment = -4108

[void]$ws.Range("A8:B8").Merge()
[void]$ws.Range("A12:B12").Merge()

[void]$ws.Range("A12:B12").Select()
